$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H108").Value = 99696
$ws.Range("J108").Value = 99696
$ws.Range("L108").Value = 99696
$ws.Range("N108").Value = -107376

$ws.Range("H109").Value = 86997.14
$ws.Range("J109").Value = 86997.14
$ws.Range("L109").Value = 86997.14
$ws.Range("N109").Value = -89771.14

$ws.Range("H110").Value = 67980.664
$ws.Range("J110").Value = 67980.664
$ws.Range("L110").Value = 67980.664
$ws.Range("N110").Value = -76160.664

$ws.Range("H112").Value = 809.6667
$ws.Range("J112").Value = 999.94446
$ws.Range("L112").Value = 2999.83338
$ws.Range("N112").Value = -5215.83338

$ws.Range("H117").Value = 83804
$ws.Range("I117").Value = 40000
$ws.Range("K117").Value = 40000
$ws.Range("M117").Value = -35411

$ws.Range("H132").Value = 1144.1428
$ws.Range("I132").Value = 1183.091
$ws.Range("J132").Value = 801.4
$ws.Range("K132").Value = 3549.273
$ws.Range("L132").Value = 2404.2
$ws.Range("M132").Value = -1019.273
$ws.Range("N132").Value = -7464.2

$ws.Range("H134").Value = 40444
$ws.Range("J134").Value = 40444
$ws.Range("L134").Value = 40444
$ws.Range("N134").Value = -50584

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4754.246
$ws.Range("I32").Value = 1645.5862
$ws.Range("K32").Value = 1645.5862
$ws.Range("M32").Value = -1358.5862

$ws.Range("H52").Value = 37190.4
$ws.Range("J52").Value = 37190.4
$ws.Range("L52").Value = 37190.4
$ws.Range("N52").Value = -37826.4

$ws.Range("H74").Value = 50569
$ws.Range("I74").Value = 68565.664
$ws.Range("K74").Value = 68565.664
$ws.Range("M74").Value = -67691.664

$ws.Range("H77").Value = 50569
$ws.Range("I77").Value = 68565.664
$ws.Range("K77").Value = 342828.32
$ws.Range("M77").Value = -338460.32

$ws.Range("H132").Value = 1750.5405
$ws.Range("I132").Value = 1655.0294
$ws.Range("K132").Value = 4965.0882
$ws.Range("M132").Value = -2435.0882

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H50").Value = 56475.75
$ws.Range("J50").Value = 58398
$ws.Range("L50").Value = 58398
$ws.Range("N50").Value = -59546

$ws.Range("H52").Value = 99990
$ws.Range("J52").Value = 99990
$ws.Range("L52").Value = 99990
$ws.Range("N52").Value = -100516

$ws.Range("H55").Value = 35246.75
$ws.Range("J55").Value = 35246.75
$ws.Range("L55").Value = 35246.75
$ws.Range("N55").Value = -35792.75

$ws.Range("H105").Value = 26325.61
$ws.Range("J105").Value = 5199.6
$ws.Range("L105").Value = 5199.6
$ws.Range("N105").Value = -8693.6

$ws.Range("H109").Value = 77996.664
$ws.Range("J109").Value = 77996.664
$ws.Range("L109").Value = 77996.664
$ws.Range("N109").Value = -80770.664

$ws.Range("H119").Value = 99396.664
$ws.Range("J119").Value = 99396.664
$ws.Range("L119").Value = 99396.664
$ws.Range("N119").Value = -109072.664

$ws.Range("H121").Value = 99990
$ws.Range("J121").Value = 99990
$ws.Range("L121").Value = 99990
$ws.Range("N121").Value = -103484

$ws.Range("H122").Value = 71183.89
$ws.Range("J122").Value = 71183.89
$ws.Range("L122").Value = 71183.89
$ws.Range("N122").Value = -80983.89

$ws.Range("H127").Value = 60384
$ws.Range("J127").Value = 60384
$ws.Range("L127").Value = 60384
$ws.Range("N127").Value = -70304

$ws.Range("H132").Value = 30646.666
$ws.Range("J132").Value = 30646.666
$ws.Range("L132").Value = 30646.666
$ws.Range("N132").Value = -40766.666

$ws.Range("H134").Value = 4554.64
$ws.Range("J134").Value = 6318.727
$ws.Range("L134").Value = 18956.181
$ws.Range("N134").Value = -24026.181

$ws.Range("H135").Value = 50706.77
$ws.Range("J135").Value = 50706.77
$ws.Range("L135").Value = 50706.77
$ws.Range("N135").Value = -60846.77

$ws.Range("H140").Value = 43499
$ws.Range("J140").Value = 43499
$ws.Range("L140").Value = 43499
$ws.Range("N140").Value = -53859

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3093.3635
$ws.Range("I31").Value = 1937.5385
$ws.Range("K31").Value = 1937.5385
$ws.Range("M31").Value = -1642.5385

$ws.Range("H34").Value = 3093.3635
$ws.Range("I34").Value = 1937.5385
$ws.Range("K34").Value = 1937.5385
$ws.Range("M34").Value = -1735.5385

$ws.Range("H108").Value = 47572.75
$ws.Range("J108").Value = 47572.75
$ws.Range("L108").Value = 47572.75
$ws.Range("N108").Value = -55252.75

$ws.Range("H116").Value = 96965.664
$ws.Range("J116").Value = 96965.664
$ws.Range("L116").Value = 96965.664
$ws.Range("N116").Value = -106143.664

$ws.Range("H119").Value = 99999
$ws.Range("J119").Value = 99999
$ws.Range("L119").Value = 99999
$ws.Range("N119").Value = -109675

$ws.Range("H134").Value = 3451822.8
$ws.Range("I134").Value = 3971210
$ws.Range("J134").Value = 335499
$ws.Range("K134").Value = 11913630
$ws.Range("L134").Value = 1006497
$ws.Range("M134").Value = -11911095
$ws.Range("N134").Value = -1011567

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H108").Value = 49739
$ws.Range("J108").Value = 49739
$ws.Range("L108").Value = 49739
$ws.Range("N108").Value = -57419

$ws.Range("H109").Value = 35666
$ws.Range("J109").Value = 35666
$ws.Range("L109").Value = 35666
$ws.Range("N109").Value = -37746

$ws.Range("H114").Value = 74515.55499999999
$ws.Range("J114").Value = 74515.55499999999
$ws.Range("L114").Value = 74515.55499999999
$ws.Range("N114").Value = -83193.55499999999

$ws.Range("H119").Value = 58192.8
$ws.Range("J119").Value = 58192.8
$ws.Range("L119").Value = 58192.8
$ws.Range("N119").Value = -67868.8

$ws.Range("H135").Value = 45436.25
$ws.Range("J135").Value = 45436.25
$ws.Range("L135").Value = 45436.25
$ws.Range("N135").Value = -55576.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4867.375
$ws.Range("I7").Value = 4456.5713
$ws.Range("J7").Value = 5186.8887
$ws.Range("K7").Value = 4456.5713
$ws.Range("L7").Value = 5186.8887
$ws.Range("M7").Value = -4344.5713
$ws.Range("N7").Value = -5410.8887

$ws.Range("H100").Value = 4537.125
$ws.Range("I100").Value = 4537.125
$ws.Range("K100").Value = 4537.125
$ws.Range("M100").Value = -3996.125

$ws.Range("H117").Value = 84096.664
$ws.Range("J117").Value = 84096.664
$ws.Range("L117").Value = 84096.664
$ws.Range("N117").Value = -93274.664

$ws.Range("H126").Value = 4867.375
$ws.Range("I126").Value = 4456.5713
$ws.Range("J126").Value = 5186.8887
$ws.Range("K126").Value = 13369.7139
$ws.Range("L126").Value = 15560.6661
$ws.Range("M126").Value = -10899.7139
$ws.Range("N126").Value = -20500.6661

$ws.Range("H136").Value = 2363.9285
$ws.Range("I136").Value = 1701.75
$ws.Range("J136").Value = 3246.8333
$ws.Range("K136").Value = 5105.25
$ws.Range("L136").Value = 9740.499899999999
$ws.Range("M136").Value = -2555.25
$ws.Range("N136").Value = -14840.4999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H121").Value = 35282.855
$ws.Range("J121").Value = 35282.855
$ws.Range("L121").Value = 35282.855
$ws.Range("N121").Value = -38776.855
